$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

$ws.Cells.Item($row, 1).Value = 88
$ws.Cells.Item($row, 2).Value = "netherlands"
$ws.Cells.Item($row, 3).Value = "eredivisie"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45232.83333333334
$ws.Cells.Item($row, 6).Value = "Ajax"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "FC Volendam"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 1.08
$ws.Cells.Item($row, 11).Value = "24/09/2023 14:42"
$ws.Cells.Item($row, 12).Value = 1.21
$ws.Cells.Item($row, 13).Value = "02/11/2023 19:59"
$ws.Cells.Item($row, 14).Value = 14.87
$ws.Cells.Item($row, 15).Value = "24/09/2023 14:42"
$ws.Cells.Item($row, 16).Value = 8.38
$ws.Cells.Item($row, 17).Value = "02/11/2023 19:59"
$ws.Cells.Item($row, 18).Value = 19.67
$ws.Cells.Item($row, 19).Value = "24/09/2023 14:42"
$ws.Cells.Item($row, 20).Value = 11.13
$ws.Cells.Item($row, 21).Value = "02/11/2023 19:59"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/ajax-fc-volendam/4tNgFLTK/"

# Carry over the formatting from the row above (column A uses the bold
# bordered "index" style, column E uses the date/time number format) by
# copying the existing cell formats onto the new row, matching how the
# rest of the sheet is styled.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($row - 1, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
